$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one weekly price observation per row (rows 2-67, newest
# data keeps getting prepended by the publisher). This week's refresh:
# insert a brand-new row at 13 (pushing the former rows 13-67 down to
# 14-68, which also grows the used range to A1:R68) and populate the new
# row with the latest observation.
$ws.Rows.Item(13).Insert()

$ws.Cells.Item(13, 1).Value = 2
$ws.Cells.Item(13, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(13, 3).Value = "Coquimbo"
$ws.Cells.Item(13, 4).Value = 44811
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 6).Value = 100112026
$ws.Cells.Item(13, 7).Value = "Haba"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 500
$ws.Cells.Item(13, 11).Value = 5000
$ws.Cells.Item(13, 12).Value = 6000
$ws.Cells.Item(13, 13).Value = 5500
$ws.Cells.Item(13, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 16).Value = 220
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
